# Merge the split word-by-word runs in the Title, Author, and Abstract
# paragraphs back into a single run each, without changing the visible
# text. Find/Replace (with identical find/replace text) collapses all
# runs spanned by the match into one run carrying the full text.

$d = $word.ActiveDocument

# Paragraph 1 (style "Title"): "Factsheet: Laws of indices"
$titlePar = $d.Paragraphs(1)
$titlePar.Range.Find.Execute("Factsheet: Laws of indices", $false, $false, $false, $false, $false,
                              $true, 1, $false, "Factsheet: Laws of indices", 2)

# Paragraph 2 (style "Author"): "Tom Coleman"
$authorPar = $d.Paragraphs(2)
$authorPar.Range.Find.Execute("Tom Coleman", $false, $false, $false, $false, $false,
                               $true, 1, $false, "Tom Coleman", 2)

# Paragraph 4 (style "Abstract"): "A list of laws of indices."
$abstractPar = $d.Paragraphs(4)
$abstractPar.Range.Find.Execute("A list of laws of indices.", $false, $false, $false, $false, $false,
                                 $true, 1, $false, "A list of laws of indices.", 2)
